$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 178-179; this pushes the previous rows 178-247
# down to 180-249 (and extends the used range to A1:R249).
$ws.Rows("178:179").Insert()

# Populate the newly inserted row 178.
$ws.Cells.Item(178, 1).Value = 11
$ws.Cells.Item(178, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(178, 3).Value = "Bíobío"
$ws.Cells.Item(178, 4).Value = 44784
$ws.Cells.Item(178, 5).Value = 8
$ws.Cells.Item(178, 6).Value = 100112045
$ws.Cells.Item(178, 7).Value = "Zapallo"
$ws.Cells.Item(178, 8).Value = "Camote"
$ws.Cells.Item(178, 9).Value = "1a (guarda)"
$ws.Cells.Item(178, 10).Value = 600
$ws.Cells.Item(178, 11).Value = 900
$ws.Cells.Item(178, 12).Value = 950
$ws.Cells.Item(178, 13).Value = 925
$ws.Cells.Item(178, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(178, 15).Value = "Región Metropolitana"
$ws.Cells.Item(178, 16).Value = 925
$ws.Cells.Item(178, 17).Value = 1
$ws.Cells.Item(178, 18).Value = "Hortaliza"

# Populate the newly inserted row 179.
$ws.Cells.Item(179, 1).Value = 11
$ws.Cells.Item(179, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(179, 3).Value = "Bíobío"
$ws.Cells.Item(179, 4).Value = 44784
$ws.Cells.Item(179, 5).Value = 8
$ws.Cells.Item(179, 6).Value = 100112045
$ws.Cells.Item(179, 7).Value = "Zapallo"
$ws.Cells.Item(179, 8).Value = "Camote"
$ws.Cells.Item(179, 9).Value = "2a (guarda)"
$ws.Cells.Item(179, 10).Value = 300
$ws.Cells.Item(179, 11).Value = 850
$ws.Cells.Item(179, 12).Value = 850
$ws.Cells.Item(179, 13).Value = 850
$ws.Cells.Item(179, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(179, 15).Value = "Región Metropolitana"
$ws.Cells.Item(179, 16).Value = 850
$ws.Cells.Item(179, 17).Value = 1
$ws.Cells.Item(179, 18).Value = "Hortaliza"
